$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '58.743.47'
$ws.Range("E2").Value = '  +0.43%  '
$ws.Range("D3").Value = '2.559.61'
$ws.Range("E3").Value = '  -0.77%  '
$ws.Range("E4").Value = '  +0.14%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '560.95'
$ws.Range("D5").NumberFormat = "General"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +3.71%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '141.58'
$ws.Range("D6").NumberFormat = "General"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.69%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("D7").NumberFormat = "General"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.17%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.589'
$ws.Range("D8").NumberFormat = "General"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.59%  '
$ws.Range("D9").Value = '2.559.74'
$ws.Range("E9").Value = '  -1.00%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.60'
$ws.Range("D10").NumberFormat = "General"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -2.39%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.101'
$ws.Range("D11").NumberFormat = "General"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.18%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.151'
$ws.Range("D12").NumberFormat = "General"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +9.12%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.336'
$ws.Range("D13").NumberFormat = "General"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.03%  '
$ws.Range("D14").Value = '3.025.99'
$ws.Range("E14").Value = '  -0.16%  '
$ws.Range("D15").Value = '58.898.19'
$ws.Range("E15").Value = '  +0.77%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '21.50'
$ws.Range("D16").NumberFormat = "General"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +4.39%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000134'
$ws.Range("D17").NumberFormat = "General"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.92%  '
$ws.Range("D18").Value = '2.583.76'
$ws.Range("E18").Value = '  -3.18%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.46'
$ws.Range("D19").NumberFormat = "General"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.55%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '332.61'
$ws.Range("D20").NumberFormat = "General"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.74%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.06'
$ws.Range("D21").NumberFormat = "General"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.14%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.13'
$ws.Range("D22").NumberFormat = "General"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.53%  '
$ws.Range("E23").Value = '  -0.02%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '64.39'
$ws.Range("D24").NumberFormat = "General"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -3.06%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.441'
$ws.Range("D25").NumberFormat = "General"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +4.39%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.00'
$ws.Range("D26").NumberFormat = "General"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.63%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.160'
$ws.Range("D27").NumberFormat = "General"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.04%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.14'
$ws.Range("D28").NumberFormat = "General"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.06%  '
$ws.Range("D29").Value = '0.0₃0771'
$ws.Range("E29").Value = '  +4.53%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.00'
$ws.Range("D30").NumberFormat = "General"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.10%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.67'
$ws.Range("D31").NumberFormat = "General"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.66%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '159.04'
$ws.Range("D32").NumberFormat = "General"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +3.60%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.96'
$ws.Range("D33").NumberFormat = "General"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.41%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '18.74'
$ws.Range("D34").NumberFormat = "General"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.17%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.96'
$ws.Range("D35").NumberFormat = "General"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.33%  '
$ws.Range("B36").Value = 'SuiNetwork'
$ws.Range("C36").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.870'
$ws.Range("D36").NumberFormat = "General"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +2.30%  '
$ws.Range("B37").Value = 'Fetch.AI'
$ws.Range("C37").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.868'
$ws.Range("D37").NumberFormat = "General"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +5.95%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.11'
$ws.Range("D38").NumberFormat = "General"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.89%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '36.74'
$ws.Range("D39").NumberFormat = "General"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.87%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.46'
$ws.Range("D40").NumberFormat = "General"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +3.09%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '288.59'
$ws.Range("D41").NumberFormat = "General"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +3.72%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.57'
$ws.Range("D42").NumberFormat = "General"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.62%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.00'
$ws.Range("D43").NumberFormat = "General"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.11%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0967'
$ws.Range("D44").NumberFormat = "General"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +2.69%  '
$ws.Range("B45").Value = 'WhiteBITCoin'
$ws.Range("C45").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '10.61'
$ws.Range("D45").NumberFormat = "General"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.18%  '
$ws.Range("B46").Value = 'Mantle'
$ws.Range("C46").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.585'
$ws.Range("D46").NumberFormat = "General"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.72%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0531'
$ws.Range("D47").NumberFormat = "General"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.17%  '
$ws.Range("B48").Value = 'Aave'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '123.88'
$ws.Range("D48").NumberFormat = "General"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +13.65%  '
$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '18.75'
$ws.Range("D49").NumberFormat = "General"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.25%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0228'
$ws.Range("D50").NumberFormat = "General"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.72%  '
$ws.Range("D51").Value = '1.919.03'
$ws.Range("E51").Value = '  +0.80%  '
